$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) retains text formatting so values like "1.00"
# are not coerced into numbers by Excel, matching the original inlineStr data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '57.709.67'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '3.064.35'
$ws.Range('E3').Value = '  +1.81%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').Value = '518.38'
$ws.Range('E5').Value = '  +1.47%  '
$ws.Range('D6').Value = '140.58'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').Value = '0.435'
$ws.Range('E8').Value = '  -0.46%  '
$ws.Range('D9').Value = '7.31'
$ws.Range('E9').Value = '  -3.04%  '
$ws.Range('D10').Value = '0.110'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D11').Value = '0.372'
$ws.Range('E11').Value = '  +1.89%  '
$ws.Range('D12').Value = '3.578.62'
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('E13').Value = '  -3.02%  '
$ws.Range('D14').Value = '26.97'
$ws.Range('E14').Value = '  +1.95%  '
$ws.Range('D15').Value = '0.0000168'
$ws.Range('E15').Value = '  +2.41%  '
$ws.Range('D16').Value = '57.706.21'
$ws.Range('E16').Value = '  +0.27%  '
$ws.Range('D17').Value = '6.23'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('D18').Value = '3.069.21'
$ws.Range('E18').Value = '  +1.97%  '
$ws.Range('D19').Value = '13.41'
$ws.Range('E19').Value = '  +4.48%  '
$ws.Range('D20').Value = '8.22'
$ws.Range('E20').Value = '  +3.05%  '
$ws.Range('D21').Value = '331.04'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('E22').Value = '  +0.28%  '
$ws.Range('D23').Value = '0.510'
$ws.Range('E23').Value = '  +2.45%  '
$ws.Range('D24').Value = '65.46'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('D25').Value = '3.177.24'
$ws.Range('E25').Value = '  +1.32%  '
$ws.Range('E26').Value = '  -1.89%  '
$ws.Range('D27').Value = '0.993'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').Value = '0.0₃0906'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('D29').Value = '6.69'
$ws.Range('E29').Value = '  -1.12%  '
$ws.Range('D30').Value = '7.29'
$ws.Range('E30').Value = '  -1.05%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '1.82'
$ws.Range('E31').Value = '  +0.36%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').Value = '1.22'
$ws.Range('E32').Value = '  +2.22%  '
$ws.Range('D33').Value = '20.91'
$ws.Range('E33').Value = '  +1.60%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '154.52'
$ws.Range('E34').Value = '  +0.11%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '4.64'
$ws.Range('E35').Value = '  -1.55%  '
$ws.Range('D36').Value = '5.90'
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('D37').Value = '25.68'
$ws.Range('E37').Value = '  +5.07%  '
$ws.Range('D38').Value = '1.28'
$ws.Range('E38').Value = '  +0.11%  '
$ws.Range('D39').Value = '0.0679'
$ws.Range('E39').Value = '  +0.54%  '
$ws.Range('D40').Value = '37.17'
$ws.Range('E40').Value = '  -0.47%  '
$ws.Range('D41').Value = '3.89'
$ws.Range('E41').Value = '  +1.27%  '
$ws.Range('B42').Value = 'FirstDigitalUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.45%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.669'
$ws.Range('E43').Value = '  +2.89%  '
$ws.Range('D44').Value = '1.40'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').Value = '2.209.71'
$ws.Range('E45').Value = '  -1.49%  '
$ws.Range('D46').Value = '6.13'
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('D47').Value = '0.961'
$ws.Range('E47').Value = '  -2.52%  '
$ws.Range('E48').Value = '  +2.57%  '
$ws.Range('D49').Value = '20.03'
$ws.Range('E49').Value = '  +3.43%  '
$ws.Range('B50').Value = 'Notcoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/2L2Y4ghjj+notcoin-not'
$ws.Range('D50').Value = '0.0174'
$ws.Range('E50').Value = '  +13.81%  '
$ws.Range('B51').Value = 'dogwifhat'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D51').Value = '1.77'
$ws.Range('E51').Value = '  -4.83%  '
